# Case_4_74 / res_bus / vm_pu.xlsx
# "case with 380 kV done": slack-bus setpoint (col B) moves from 1.05 to 1.02 p.u.
# and the dependent bus-voltage results (cols C-F, I-N) are refreshed with the
# recomputed power-flow solution. Column G stays 1 and column H stays blank.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2

# Each inner array holds, in order: B,C,D,E,F,I,J,K,L,M,N for one sheet row.
$rowData = @(
    @(1.02,1.054703397297098,1.06066765044274,1.058387495175654,1.069392940101595,1.049390810923261,1.05971410419369,1.063393805474529,1.061119867658763,1.072095568152594,1.023586948131093),
    @(1.02,1.055840672272789,1.061587128688237,1.059475010599624,1.070444297958106,1.049741470389257,1.060501856882689,1.064127489012212,1.062020712803779,1.072962504252294,1.023858171336445),
    @(1.02,1.056576607756947,1.062182114724941,1.060179077789721,1.071124926992027,1.049967219603939,1.061011058796739,1.064601622144134,1.062603395573644,1.073523193493495,1.024033287165428),
    @(1.02,1.056886006064761,1.062432252186419,1.060475157025801,1.071411142627186,1.050061849143899,1.061225001466841,1.064800802046046,1.062848302028143,1.073758841252598,1.024106813955276),
    @(1.02,1.056937956073672,1.062474251680534,1.060524875283962,1.071459204139533,1.050077721711652,1.061260916013222,1.064834236691069,1.06288941977823,1.073798403683076,1.024119154036302),
    @(1.02,1.056580741907859,1.062185457052947,1.060183033663297,1.071128751108011,1.049968485130604,1.061013918003163,1.064604284167883,1.062606668235035,1.073526342489494,1.024034269995164),
    @(1.02,1.055087736404329,1.060978387859802,1.05875494928257,1.069748183125584,1.04950955629434,1.059980438145284,1.06364188349841,1.061424359117806,1.072388610595822,1.023678688710707),
    @(1.02,1.052457157422065,1.058851545146872,1.056241315771014,1.067317968153501,1.048692045813814,1.058155266942254,1.061941336934925,1.059339247430486,1.070381656842992,1.023049171825996),
    @(1.02,1.050703575389402,1.057433756832959,1.054567440897335,1.065699508217251,1.048141101382746,1.056935743699301,1.06080448333756,1.057947985989279,1.069042245252828,1.022627520240493),
    @(1.02,1.049944274309114,1.056819861993906,1.053843074369796,1.064999092763761,1.047901125873262,1.056407020773768,1.06031146086694,1.057345267069496,1.068461919210668,1.0224444722343),
    @(1.02,1.04966223656468,1.05659183635658,1.053574076693299,1.064738985200712,1.047811775749946,1.056210529620004,1.06012821616405,1.057121345725292,1.068246306959083,1.022376409378286),
    @(1.02,1.049722734616297,1.056640748533205,1.05363177470857,1.064794776541694,1.047830951276428,1.056252682170016,1.060167527949645,1.057169379646218,1.068292558938869,1.022391012291023),
    @(1.02,1.049920960978479,1.0568010132796,1.05382083762991,1.064977591012328,1.047893744501305,1.056390780795247,1.060296316134821,1.057326758576697,1.068444097724948,1.022438847579461),
    @(1.02,1.050043094809551,1.056899757967257,1.053937334037131,1.065090236703833,1.047932405309104,1.056475854659183,1.060375651640358,1.057423718988607,1.068537458693869,1.022468311112604),
    @(1.02,1.050753967837633,1.057474499341771,1.054615523771529,1.065746000630371,1.048156997976106,1.056970819317372,1.060837187622569,1.057987980223215,1.069080752127609,1.022639658625271),
    @(1.02,1.05119988191004,1.057835023728165,1.055041049565157,1.066157447695212,1.048297500558045,1.057281120120796,1.061126493886722,1.058341847376219,1.069421451198984,1.022747014501563),
    @(1.02,1.051459977238325,1.058045313149309,1.055289293319249,1.066397475308633,1.048379317026353,1.057462049560104,1.061295168289779,1.05854822381805,1.069620141271135,1.02280958801475),
    @(1.02,1.051548663264239,1.058117016695665,1.055373945126843,1.066479324895854,1.048407191214902,1.057523731023998,1.061352669525187,1.058618588072713,1.069687883709018,1.022830916274564),
    @(1.02,1.051152039435442,1.057796342705192,1.054995390359947,1.066113299443022,1.048282440061737,1.057247834381153,1.061095461620969,1.058303883709159,1.069384900925753,1.022735500923455),
    @(1.02,1.049862588252085,1.056753819246812,1.053765161541947,1.064923755111745,1.047875259330696,1.056350116954871,1.060258394380252,1.057280415643036,1.06839947480488,1.022424763230704),
    @(1.02,1.049051861304588,1.056098356183029,1.052992039579868,1.064176176287515,1.047618018861076,1.055785108006965,1.059731436017598,1.056636661723518,1.06777958957848,1.022228981173254),
    @(1.02,1.049481643206541,1.056445828438488,1.053401851141713,1.064572450364004,1.047754503555704,1.056084684951331,1.060010849408223,1.056977952600142,1.06810823187445,1.022332807725904),
    @(1.02,1.051173657386292,1.057813820989201,1.055016021660122,1.066133248021759,1.048289245672805,1.057262874970943,1.061109483991602,1.058321037943323,1.069401416525575,1.022740703553883),
    @(1.02,1.053137196850874,1.059401365714768,1.056890816996458,1.067945939358413,1.048904437810381,1.058627599282065,1.062381524233745,1.059878506347578,1.070900755717884,1.023212264984273)
)

$lastRow = $firstRow + $rowData.Count - 1
$bf  = New-Object "object[,]" $rowData.Count,5
$inn = New-Object "object[,]" $rowData.Count,6
for ($i = 0; $i -lt $rowData.Count; $i++) {
    $r = $rowData[$i]
    for ($j = 0; $j -lt 5; $j++) { $bf[$i, $j] = $r[$j] }
    for ($j = 0; $j -lt 6; $j++) { $inn[$i, $j] = $r[5 + $j] }
}

$ws.Range("B" + $firstRow + ":F" + $lastRow).Value = $bf
$ws.Range("I" + $firstRow + ":N" + $lastRow).Value = $inn
